$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 35.916668
$ws.Cells.Item(6, 9).Value = 35.916668
$ws.Cells.Item(6, 11).Value = 107.750004
$ws.Cells.Item(6, 13).Value = 4.249995999999996
$ws.Cells.Item(7, 8).Value = 13950
$ws.Cells.Item(7, 10).Value = 20000
$ws.Cells.Item(7, 12).Value = 20000
$ws.Cells.Item(7, 14).Value = -20224
$ws.Cells.Item(9, 8).Value = 260.5
$ws.Cells.Item(9, 10).Value = 550
$ws.Cells.Item(9, 12).Value = 550
$ws.Cells.Item(9, 14).Value = -888
$ws.Cells.Item(12, 8).Value = 1267.3334
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 1267.3334
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 1267.3334
$ws.Cells.Item(12, 13).Value = $null
$ws.Cells.Item(12, 14).Value = -1607.3334
$ws.Cells.Item(14, 8).Value = 13950
$ws.Cells.Item(14, 10).Value = 20000
$ws.Cells.Item(14, 12).Value = 20000
$ws.Cells.Item(14, 14).Value = -20382
$ws.Cells.Item(28, 8).Value = 4834.706
$ws.Cells.Item(28, 9).Value = 1602.25
$ws.Cells.Item(28, 10).Value = 7708
$ws.Cells.Item(28, 11).Value = 1602.25
$ws.Cells.Item(28, 12).Value = 7708
$ws.Cells.Item(28, 13).Value = -1117.25
$ws.Cells.Item(28, 14).Value = -8678
$ws.Cells.Item(40, 8).Value = 5456.5386
$ws.Cells.Item(40, 9).Value = 4485
$ws.Cells.Item(40, 10).Value = 6289.2856
$ws.Cells.Item(40, 11).Value = 4485
$ws.Cells.Item(40, 12).Value = 6289.2856
$ws.Cells.Item(40, 13).Value = -4310
$ws.Cells.Item(40, 14).Value = -6639.2856
$ws.Cells.Item(51, 8).Value = 3000
$ws.Cells.Item(51, 9).Value = 3000
$ws.Cells.Item(51, 10).Value = 3000
$ws.Cells.Item(51, 11).Value = 3000
$ws.Cells.Item(51, 12).Value = 3000
$ws.Cells.Item(51, 13).Value = -2516
$ws.Cells.Item(51, 14).Value = -3968
$ws.Cells.Item(58, 8).Value = 2236
$ws.Cells.Item(58, 9).Value = 1715
$ws.Cells.Item(58, 10).Value = 2583.3333
$ws.Cells.Item(58, 11).Value = 5145
$ws.Cells.Item(58, 12).Value = 7749.999899999999
$ws.Cells.Item(58, 13).Value = -4995
$ws.Cells.Item(58, 14).Value = -8049.999899999999
$ws.Cells.Item(69, 8).Value = 7632.1816
$ws.Cells.Item(69, 10).Value = 7757.524
$ws.Cells.Item(69, 12).Value = 23272.572
$ws.Cells.Item(69, 14).Value = -25020.572
$ws.Cells.Item(70, 8).Value = 7947.125
$ws.Cells.Item(70, 9).Value = 3875
$ws.Cells.Item(70, 10).Value = 8528.857
$ws.Cells.Item(70, 11).Value = 11625
$ws.Cells.Item(70, 12).Value = 25586.571
$ws.Cells.Item(70, 13).Value = -11355
$ws.Cells.Item(70, 14).Value = -26126.571
$ws.Cells.Item(72, 8).Value = 7632.1816
$ws.Cells.Item(72, 10).Value = 7757.524
$ws.Cells.Item(72, 12).Value = 69817.716
$ws.Cells.Item(72, 14).Value = -78553.716
$ws.Cells.Item(73, 8).Value = 7947.125
$ws.Cells.Item(73, 9).Value = 3875
$ws.Cells.Item(73, 10).Value = 8528.857
$ws.Cells.Item(73, 11).Value = 11625
$ws.Cells.Item(73, 12).Value = 25586.571
$ws.Cells.Item(73, 13).Value = -10689
$ws.Cells.Item(73, 14).Value = -27458.571
$ws.Cells.Item(74, 8).Value = 7399
$ws.Cells.Item(74, 9).Value = 7399
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 7399
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = -6463
$ws.Cells.Item(74, 14).Value = $null
$ws.Cells.Item(77, 8).Value = 7399
$ws.Cells.Item(77, 9).Value = 7399
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 36995
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = -32315
$ws.Cells.Item(77, 14).Value = $null
$ws.Cells.Item(88, 8).Value = 1488
$ws.Cells.Item(88, 9).Value = 1391
$ws.Cells.Item(88, 10).Value = 1649.6666
$ws.Cells.Item(88, 11).Value = 1391
$ws.Cells.Item(88, 12).Value = 1649.6666
$ws.Cells.Item(88, 13).Value = -985
$ws.Cells.Item(88, 14).Value = -2461.6666
$ws.Cells.Item(91, 8).Value = 1488
$ws.Cells.Item(91, 9).Value = 1391
$ws.Cells.Item(91, 10).Value = 1649.6666
$ws.Cells.Item(91, 11).Value = 1391
$ws.Cells.Item(91, 12).Value = 1649.6666
$ws.Cells.Item(91, 13).Value = 13
$ws.Cells.Item(91, 14).Value = -4457.6666
$ws.Cells.Item(95, 8).Value = 26399.4
$ws.Cells.Item(95, 10).Value = 26399.4
$ws.Cells.Item(95, 12).Value = 26399.4
$ws.Cells.Item(95, 14).Value = -31891.4
$ws.Cells.Item(98, 8).Value = 568.4167
$ws.Cells.Item(98, 9).Value = 483.72726
$ws.Cells.Item(98, 11).Value = 483.72726
$ws.Cells.Item(98, 13).Value = 1014.27274
$ws.Cells.Item(100, 8).Value = 2311.8462
$ws.Cells.Item(100, 9).Value = 2431.75
$ws.Cells.Item(100, 11).Value = 2431.75
$ws.Cells.Item(100, 13).Value = -1890.75
$ws.Cells.Item(105, 8).Value = 16063.111
$ws.Cells.Item(105, 10).Value = 16063.111
$ws.Cells.Item(105, 12).Value = 16063.111
$ws.Cells.Item(105, 14).Value = -23051.111
$ws.Cells.Item(107, 8).Value = 194.08333
$ws.Cells.Item(107, 9).Value = 194.33333
$ws.Cells.Item(107, 10).Value = 193.33333
$ws.Cells.Item(107, 11).Value = 194.33333
$ws.Cells.Item(107, 12).Value = 193.33333
$ws.Cells.Item(107, 13).Value = 1725.66667
$ws.Cells.Item(107, 14).Value = -4033.33333
$ws.Cells.Item(112, 8).Value = 2631.7693
$ws.Cells.Item(112, 10).Value = 2523.6667
$ws.Cells.Item(112, 12).Value = 7571.000100000001
$ws.Cells.Item(112, 14).Value = -9787.000100000001
$ws.Cells.Item(113, 8).Value = 7999.6665
$ws.Cells.Item(113, 10).Value = 1999.5
$ws.Cells.Item(113, 12).Value = 1999.5
$ws.Cells.Item(113, 14).Value = -8507.5
$ws.Cells.Item(122, 8).Value = 568.4167
$ws.Cells.Item(122, 9).Value = 483.72726
$ws.Cells.Item(122, 11).Value = 1451.18178
$ws.Cells.Item(122, 13).Value = 998.8182200000001
$ws.Cells.Item(127, 8).Value = 818.75
$ws.Cells.Item(127, 9).Value = 758.3333
$ws.Cells.Item(127, 10).Value = 1000
$ws.Cells.Item(127, 11).Value = 2274.9999
$ws.Cells.Item(127, 12).Value = 3000
$ws.Cells.Item(127, 13).Value = 2685.0001
$ws.Cells.Item(127, 14).Value = -12920
$ws.Cells.Item(129, 8).Value = 883.5
$ws.Cells.Item(129, 9).Value = 581.7143
$ws.Cells.Item(129, 10).Value = 2996
$ws.Cells.Item(129, 11).Value = 1745.1429
$ws.Cells.Item(129, 12).Value = 8988
$ws.Cells.Item(129, 13).Value = 3254.8571
$ws.Cells.Item(129, 14).Value = -18988
$ws.Cells.Item(130, 8).Value = 100000
$ws.Cells.Item(130, 10).Value = 100000
$ws.Cells.Item(130, 12).Value = 100000
$ws.Cells.Item(130, 14).Value = -110040
$ws.Cells.Item(132, 8).Value = 13508.529
$ws.Cells.Item(132, 9).Value = 14277
$ws.Cells.Item(132, 10).Value = 7745
$ws.Cells.Item(132, 11).Value = 42831
$ws.Cells.Item(132, 12).Value = 23235
$ws.Cells.Item(132, 13).Value = -40301
$ws.Cells.Item(132, 14).Value = -28295
$ws.Cells.Item(135, 8).Value = 1771.4286
$ws.Cells.Item(135, 9).Value = 1560.6666
$ws.Cells.Item(135, 10).Value = 3036
$ws.Cells.Item(135, 11).Value = 14045.9994
$ws.Cells.Item(135, 12).Value = 27324
$ws.Cells.Item(135, 13).Value = -11510.9994
$ws.Cells.Item(135, 14).Value = -32394
$ws.Cells.Item(137, 8).Value = 2417.5625
$ws.Cells.Item(137, 9).Value = 1288.25
$ws.Cells.Item(137, 10).Value = 2794
$ws.Cells.Item(137, 11).Value = 3864.75
$ws.Cells.Item(137, 12).Value = 8382
$ws.Cells.Item(137, 13).Value = -1314.75
$ws.Cells.Item(137, 14).Value = -13482
$ws.Cells.Item(138, 8).Value = 4143.48
$ws.Cells.Item(138, 9).Value = 1136
$ws.Cells.Item(138, 10).Value = 5835.1875
$ws.Cells.Item(138, 11).Value = 3408
$ws.Cells.Item(138, 12).Value = 17505.5625
$ws.Cells.Item(138, 13).Value = 1732
$ws.Cells.Item(138, 14).Value = -27785.5625
$ws.Cells.Item(141, 8).Value = 2152.6
$ws.Cells.Item(141, 9).Value = 2003.375
$ws.Cells.Item(141, 10).Value = 2749.5
$ws.Cells.Item(141, 11).Value = 6010.125
$ws.Cells.Item(141, 12).Value = 8248.5
$ws.Cells.Item(141, 13).Value = -830.125
$ws.Cells.Item(141, 14).Value = -18608.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 4844.5713
$ws.Cells.Item(2, 9).Value = 955
$ws.Cells.Item(2, 11).Value = 955
$ws.Cells.Item(2, 13).Value = -842
$ws.Cells.Item(32, 8).Value = 5001500
$ws.Cells.Item(32, 9).Value = 1578.7894
$ws.Cells.Item(32, 11).Value = 1578.7894
$ws.Cells.Item(32, 13).Value = -1291.7894
$ws.Cells.Item(43, 8).Value = 4314331
$ws.Cells.Item(43, 10).Value = 42496
$ws.Cells.Item(43, 12).Value = 42496
$ws.Cells.Item(43, 14).Value = -43122
$ws.Cells.Item(61, 8).Value = 2587.375
$ws.Cells.Item(61, 9).Value = 2642.7144
$ws.Cells.Item(61, 10).Value = 2200
$ws.Cells.Item(61, 11).Value = 2642.7144
$ws.Cells.Item(61, 12).Value = 2200
$ws.Cells.Item(61, 13).Value = -2430.7144
$ws.Cells.Item(61, 14).Value = -2624
$ws.Cells.Item(74, 8).Value = 6313
$ws.Cells.Item(74, 9).Value = 5910.778
$ws.Cells.Item(74, 11).Value = 5910.778
$ws.Cells.Item(74, 13).Value = -5036.778
$ws.Cells.Item(77, 8).Value = 6313
$ws.Cells.Item(77, 9).Value = 5910.778
$ws.Cells.Item(77, 11).Value = 29553.89
$ws.Cells.Item(77, 13).Value = -25185.89
$ws.Cells.Item(88, 8).Value = 1679.4
$ws.Cells.Item(88, 9).Value = 2124.8333
$ws.Cells.Item(88, 10).Value = 1011.25
$ws.Cells.Item(88, 11).Value = 2124.8333
$ws.Cells.Item(88, 12).Value = 1011.25
$ws.Cells.Item(88, 13).Value = -1718.8333
$ws.Cells.Item(88, 14).Value = -1823.25
$ws.Cells.Item(91, 8).Value = 1679.4
$ws.Cells.Item(91, 9).Value = 2124.8333
$ws.Cells.Item(91, 10).Value = 1011.25
$ws.Cells.Item(91, 11).Value = 2124.8333
$ws.Cells.Item(91, 12).Value = 1011.25
$ws.Cells.Item(91, 13).Value = -720.8332999999998
$ws.Cells.Item(91, 14).Value = -3819.25
$ws.Cells.Item(102, 8).Value = 7842
$ws.Cells.Item(102, 9).Value = 4947
$ws.Cells.Item(102, 11).Value = 4947
$ws.Cells.Item(102, 13).Value = -3325
$ws.Cells.Item(110, 8).Value = 4076.8333
$ws.Cells.Item(110, 10).Value = 6725
$ws.Cells.Item(110, 12).Value = 6725
$ws.Cells.Item(110, 14).Value = -10815
$ws.Cells.Item(116, 8).Value = 4844.5713
$ws.Cells.Item(116, 9).Value = 955
$ws.Cells.Item(116, 11).Value = 955
$ws.Cells.Item(116, 13).Value = 1339
$ws.Cells.Item(132, 8).Value = 1727.4546
$ws.Cells.Item(132, 9).Value = 1811.3334
$ws.Cells.Item(132, 10).Value = 1350
$ws.Cells.Item(132, 11).Value = 5434.0002
$ws.Cells.Item(132, 12).Value = 4050
$ws.Cells.Item(132, 13).Value = -2904.0002
$ws.Cells.Item(132, 14).Value = -9110
$ws.Cells.Item(136, 8).Value = 2587.375
$ws.Cells.Item(136, 9).Value = 2642.7144
$ws.Cells.Item(136, 10).Value = 2200
$ws.Cells.Item(136, 11).Value = 7928.1432
$ws.Cells.Item(136, 12).Value = 6600
$ws.Cells.Item(136, 13).Value = -5378.1432
$ws.Cells.Item(136, 14).Value = -11700

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 4844.5713
$ws.Cells.Item(3, 9).Value = 955
$ws.Cells.Item(3, 11).Value = 955
$ws.Cells.Item(3, 13).Value = -841
$ws.Cells.Item(80, 8).Value = 981.7
$ws.Cells.Item(80, 9).Value = 620
$ws.Cells.Item(80, 10).Value = 1524.25
$ws.Cells.Item(80, 11).Value = 620
$ws.Cells.Item(80, 12).Value = 1524.25
$ws.Cells.Item(80, 13).Value = 378
$ws.Cells.Item(80, 14).Value = -3520.25
$ws.Cells.Item(83, 8).Value = 981.7
$ws.Cells.Item(83, 9).Value = 620
$ws.Cells.Item(83, 10).Value = 1524.25
$ws.Cells.Item(83, 11).Value = 3100
$ws.Cells.Item(83, 12).Value = 7621.25
$ws.Cells.Item(83, 13).Value = 1892
$ws.Cells.Item(83, 14).Value = -17605.25
$ws.Cells.Item(86, 8).Value = 3749
$ws.Cells.Item(86, 9).Value = 2498.8
$ws.Cells.Item(86, 10).Value = 10000
$ws.Cells.Item(86, 11).Value = 2498.8
$ws.Cells.Item(86, 12).Value = 10000
$ws.Cells.Item(86, 13).Value = -1375.8
$ws.Cells.Item(86, 14).Value = -12246
$ws.Cells.Item(89, 8).Value = 3749
$ws.Cells.Item(89, 9).Value = 2498.8
$ws.Cells.Item(89, 10).Value = 10000
$ws.Cells.Item(89, 11).Value = 12494
$ws.Cells.Item(89, 12).Value = 50000
$ws.Cells.Item(89, 13).Value = -6878
$ws.Cells.Item(89, 14).Value = -61232
$ws.Cells.Item(94, 8).Value = 365
$ws.Cells.Item(94, 9).Value = 331.25
$ws.Cells.Item(94, 11).Value = 331.25
$ws.Cells.Item(94, 13).Value = 119.75
$ws.Cells.Item(134, 8).Value = 8332.875
$ws.Cells.Item(134, 9).Value = 1009
$ws.Cells.Item(134, 11).Value = 3027
$ws.Cells.Item(134, 13).Value = -492

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 65.78261000000001
$ws.Cells.Item(7, 9).Value = 49
$ws.Cells.Item(7, 11).Value = 49
$ws.Cells.Item(7, 13).Value = 64
$ws.Cells.Item(16, 8).Value = 1454.25
$ws.Cells.Item(16, 9).Value = 1454.25
$ws.Cells.Item(16, 11).Value = 1454.25
$ws.Cells.Item(16, 13).Value = -1167.25
$ws.Cells.Item(22, 8).Value = 4222.5
$ws.Cells.Item(22, 9).Value = 4001
$ws.Cells.Item(22, 11).Value = 4001
$ws.Cells.Item(22, 13).Value = -3651
$ws.Cells.Item(31, 8).Value = 5580.48
$ws.Cells.Item(31, 9).Value = 2390.2222
$ws.Cells.Item(31, 10).Value = 7375
$ws.Cells.Item(31, 11).Value = 2390.2222
$ws.Cells.Item(31, 12).Value = 7375
$ws.Cells.Item(31, 13).Value = -2095.2222
$ws.Cells.Item(31, 14).Value = -7965
$ws.Cells.Item(33, 8).Value = 974.5
$ws.Cells.Item(33, 9).Value = 974.5
$ws.Cells.Item(33, 11).Value = 974.5
$ws.Cells.Item(33, 13).Value = -595.5
$ws.Cells.Item(34, 8).Value = 5580.48
$ws.Cells.Item(34, 9).Value = 2390.2222
$ws.Cells.Item(34, 10).Value = 7375
$ws.Cells.Item(34, 11).Value = 2390.2222
$ws.Cells.Item(34, 12).Value = 7375
$ws.Cells.Item(34, 13).Value = -2188.2222
$ws.Cells.Item(34, 14).Value = -7779
$ws.Cells.Item(58, 8).Value = 5468.2
$ws.Cells.Item(58, 9).Value = 3783
$ws.Cells.Item(58, 11).Value = 3783
$ws.Cells.Item(58, 13).Value = -3580
$ws.Cells.Item(62, 8).Value = 2935
$ws.Cells.Item(62, 10).Value = 2935
$ws.Cells.Item(62, 12).Value = 2935
$ws.Cells.Item(62, 14).Value = -4183
$ws.Cells.Item(65, 8).Value = 2935
$ws.Cells.Item(65, 10).Value = 2935
$ws.Cells.Item(65, 12).Value = 14675
$ws.Cells.Item(65, 14).Value = -20915
$ws.Cells.Item(70, 8).Value = 25500
$ws.Cells.Item(70, 10).Value = 25500
$ws.Cells.Item(70, 12).Value = 25500
$ws.Cells.Item(70, 14).Value = -26130
$ws.Cells.Item(73, 8).Value = 25500
$ws.Cells.Item(73, 10).Value = 25500
$ws.Cells.Item(73, 12).Value = 25500
$ws.Cells.Item(73, 14).Value = -27684
$ws.Cells.Item(93, 8).Value = 14571
$ws.Cells.Item(93, 9).Value = 10601.75
$ws.Cells.Item(93, 11).Value = 10601.75
$ws.Cells.Item(93, 13).Value = -8729.75
$ws.Cells.Item(112, 8).Value = 39851
$ws.Cells.Item(112, 10).Value = 39851
$ws.Cells.Item(112, 12).Value = 39851
$ws.Cells.Item(112, 14).Value = -42805
$ws.Cells.Item(113, 8).Value = 1454.25
$ws.Cells.Item(113, 9).Value = 1454.25
$ws.Cells.Item(113, 11).Value = 1454.25
$ws.Cells.Item(113, 13).Value = 715.75
$ws.Cells.Item(132, 8).Value = 3515
$ws.Cells.Item(132, 10).Value = 6000
$ws.Cells.Item(132, 12).Value = 18000
$ws.Cells.Item(132, 14).Value = -23060
$ws.Cells.Item(136, 8).Value = 5468.2
$ws.Cells.Item(136, 9).Value = 3783
$ws.Cells.Item(136, 11).Value = 11349
$ws.Cells.Item(136, 13).Value = -8799

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 24.117647
$ws.Cells.Item(2, 9).Value = 21.666666
$ws.Cells.Item(2, 10).Value = 26.875
$ws.Cells.Item(2, 11).Value = 129.999996
$ws.Cells.Item(2, 12).Value = 161.25
$ws.Cells.Item(2, 13).Value = -16.99999600000001
$ws.Cells.Item(2, 14).Value = -387.25
$ws.Cells.Item(7, 8).Value = 122.8
$ws.Cells.Item(7, 9).Value = 191
$ws.Cells.Item(7, 10).Value = 77.333336
$ws.Cells.Item(7, 11).Value = 573
$ws.Cells.Item(7, 12).Value = 232.000008
$ws.Cells.Item(7, 13).Value = -461
$ws.Cells.Item(7, 14).Value = -456.000008
$ws.Cells.Item(12, 8).Value = 158.625
$ws.Cells.Item(12, 10).Value = 188.66667
$ws.Cells.Item(12, 12).Value = 566.00001
$ws.Cells.Item(12, 14).Value = -912.00001
$ws.Cells.Item(50, 8).Value = 1000
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 14).Value = $null
$ws.Cells.Item(53, 8).Value = 1000
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 14).Value = $null
$ws.Cells.Item(74, 8).Value = 9750
$ws.Cells.Item(74, 10).Value = 9750
$ws.Cells.Item(74, 12).Value = 29250
$ws.Cells.Item(74, 14).Value = -31372
$ws.Cells.Item(77, 8).Value = 9750
$ws.Cells.Item(77, 10).Value = 9750
$ws.Cells.Item(77, 12).Value = 87750
$ws.Cells.Item(77, 14).Value = -98358
$ws.Cells.Item(80, 8).Value = 4335.5654
$ws.Cells.Item(80, 10).Value = 5086.5
$ws.Cells.Item(80, 12).Value = 15259.5
$ws.Cells.Item(80, 14).Value = -17131.5
$ws.Cells.Item(83, 8).Value = 4335.5654
$ws.Cells.Item(83, 10).Value = 5086.5
$ws.Cells.Item(83, 12).Value = 45778.5
$ws.Cells.Item(83, 14).Value = -55138.5
$ws.Cells.Item(86, 8).Value = 615.8333
$ws.Cells.Item(86, 9).Value = 615.8333
$ws.Cells.Item(86, 11).Value = 1847.4999
$ws.Cells.Item(86, 13).Value = -661.4999
$ws.Cells.Item(89, 8).Value = 615.8333
$ws.Cells.Item(89, 9).Value = 615.8333
$ws.Cells.Item(89, 11).Value = 5542.4997
$ws.Cells.Item(89, 13).Value = 385.5002999999997
$ws.Cells.Item(92, 8).Value = 2417.818
$ws.Cells.Item(92, 10).Value = 3270.8572
$ws.Cells.Item(92, 12).Value = 9812.571599999999
$ws.Cells.Item(92, 14).Value = -12308.5716
$ws.Cells.Item(137, 8).Value = 5766
$ws.Cells.Item(137, 10).Value = 5766
$ws.Cells.Item(137, 12).Value = 17298
$ws.Cells.Item(137, 14).Value = -27498
$ws.Cells.Item(139, 8).Value = 3000
$ws.Cells.Item(139, 10).Value = 5000
$ws.Cells.Item(139, 12).Value = 15000
$ws.Cells.Item(139, 14).Value = -25280
$ws.Cells.Item(140, 8).Value = 2599.5715
$ws.Cells.Item(140, 9).Value = 1810.5
$ws.Cells.Item(140, 11).Value = 5431.5
$ws.Cells.Item(140, 13).Value = -251.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(38, 8).Value = 20000
$ws.Cells.Item(38, 10).Value = 20000
$ws.Cells.Item(38, 12).Value = 20000
$ws.Cells.Item(38, 14).Value = -20926
$ws.Cells.Item(80, 8).Value = 2476.25
$ws.Cells.Item(80, 9).Value = 2201.6667
$ws.Cells.Item(80, 10).Value = 3300
$ws.Cells.Item(80, 11).Value = 2201.6667
$ws.Cells.Item(80, 12).Value = 3300
$ws.Cells.Item(80, 13).Value = -1203.6667
$ws.Cells.Item(80, 14).Value = -5296
$ws.Cells.Item(83, 8).Value = 2476.25
$ws.Cells.Item(83, 9).Value = 2201.6667
$ws.Cells.Item(83, 10).Value = 3300
$ws.Cells.Item(83, 11).Value = 11008.3335
$ws.Cells.Item(83, 12).Value = 16500
$ws.Cells.Item(83, 13).Value = -6016.333500000001
$ws.Cells.Item(83, 14).Value = -26484
$ws.Cells.Item(104, 8).Value = 0
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 14).Value = $null
$ws.Cells.Item(122, 8).Value = 2053.8462
$ws.Cells.Item(122, 9).Value = 1422.8572
$ws.Cells.Item(122, 10).Value = 2790
$ws.Cells.Item(122, 11).Value = 4268.571599999999
$ws.Cells.Item(122, 12).Value = 8370
$ws.Cells.Item(122, 13).Value = -1818.571599999999
$ws.Cells.Item(122, 14).Value = -13270
$ws.Cells.Item(132, 8).Value = 2528.6667
$ws.Cells.Item(132, 9).Value = 2528.6667
$ws.Cells.Item(132, 11).Value = 7586.000100000001
$ws.Cells.Item(132, 13).Value = -5056.000100000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value = 1166.6666
$ws.Cells.Item(3, 9).Value = 1000
$ws.Cells.Item(3, 10).Value = 2000
$ws.Cells.Item(3, 11).Value = 1000
$ws.Cells.Item(3, 12).Value = 2000
$ws.Cells.Item(3, 13).Value = -888
$ws.Cells.Item(3, 14).Value = -2224
$ws.Cells.Item(7, 8).Value = 3466.6667
$ws.Cells.Item(7, 9).Value = 2700
$ws.Cells.Item(7, 10).Value = 5000
$ws.Cells.Item(7, 11).Value = 2700
$ws.Cells.Item(7, 12).Value = 5000
$ws.Cells.Item(7, 13).Value = -2588
$ws.Cells.Item(7, 14).Value = -5224
$ws.Cells.Item(15, 8).Value = 1166.6666
$ws.Cells.Item(15, 9).Value = 1000
$ws.Cells.Item(15, 10).Value = 2000
$ws.Cells.Item(15, 11).Value = 1000
$ws.Cells.Item(15, 12).Value = 2000
$ws.Cells.Item(15, 13).Value = -830
$ws.Cells.Item(15, 14).Value = -2340
$ws.Cells.Item(40, 8).Value = 3889.1785
$ws.Cells.Item(40, 9).Value = 5033
$ws.Cells.Item(40, 10).Value = 3347.3684
$ws.Cells.Item(40, 11).Value = 5033
$ws.Cells.Item(40, 12).Value = 3347.3684
$ws.Cells.Item(40, 13).Value = -4897
$ws.Cells.Item(40, 14).Value = -3619.3684
$ws.Cells.Item(46, 8).Value = 4928.1
$ws.Cells.Item(46, 9).Value = 1461.1666
$ws.Cells.Item(46, 10).Value = 6413.9287
$ws.Cells.Item(46, 11).Value = 1461.1666
$ws.Cells.Item(46, 12).Value = 6413.9287
$ws.Cells.Item(46, 13).Value = -1273.1666
$ws.Cells.Item(46, 14).Value = -6789.9287
$ws.Cells.Item(50, 8).Value = 67156
$ws.Cells.Item(50, 9).Value = 67156
$ws.Cells.Item(50, 11).Value = 67156
$ws.Cells.Item(50, 13).Value = -66519
$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 13).Value = $null
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).Value = $null
$ws.Cells.Item(82, 8).Value = 2566
$ws.Cells.Item(82, 9).Value = 415.6
$ws.Cells.Item(82, 11).Value = 415.6
$ws.Cells.Item(82, 13).Value = -54.60000000000002
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).Value = $null
$ws.Cells.Item(85, 8).Value = 2566
$ws.Cells.Item(85, 9).Value = 415.6
$ws.Cells.Item(85, 11).Value = 415.6
$ws.Cells.Item(85, 13).Value = 832.4
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).Value = $null
$ws.Cells.Item(126, 8).Value = 3466.6667
$ws.Cells.Item(126, 9).Value = 2700
$ws.Cells.Item(126, 10).Value = 5000
$ws.Cells.Item(126, 11).Value = 8100
$ws.Cells.Item(126, 12).Value = 15000
$ws.Cells.Item(126, 13).Value = -5630
$ws.Cells.Item(126, 14).Value = -19940
$ws.Cells.Item(132, 8).Value = 3659.4
$ws.Cells.Item(132, 9).Value = 3732.3333
$ws.Cells.Item(132, 11).Value = 11196.9999
$ws.Cells.Item(132, 13).Value = -8666.999899999999
$ws.Cells.Item(136, 8).Value = 5001.3335
$ws.Cells.Item(136, 9).Value = 5001.3335
$ws.Cells.Item(136, 11).Value = 15004.0005
$ws.Cells.Item(136, 13).Value = -12454.0005

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 8789.846
$ws.Cells.Item(4, 9).Value = 16733.334
$ws.Cells.Item(4, 10).Value = 1981.1428
$ws.Cells.Item(4, 11).Value = 16733.334
$ws.Cells.Item(4, 12).Value = 1981.1428
$ws.Cells.Item(4, 13).Value = -16620.334
$ws.Cells.Item(4, 14).Value = -2207.1428
$ws.Cells.Item(12, 8).Value = 1000
$ws.Cells.Item(12, 10).Value = 1000
$ws.Cells.Item(12, 12).Value = 1000
$ws.Cells.Item(12, 14).Value = -1284
$ws.Cells.Item(28, 8).Value = 19166.666
$ws.Cells.Item(28, 10).Value = 19166.666
$ws.Cells.Item(28, 12).Value = 19166.666
$ws.Cells.Item(28, 14).Value = -19862.666
$ws.Cells.Item(62, 8).Value = 9999.857
$ws.Cells.Item(62, 9).Value = 4999.5
$ws.Cells.Item(62, 10).Value = 12000
$ws.Cells.Item(62, 11).Value = 4999.5
$ws.Cells.Item(62, 12).Value = 12000
$ws.Cells.Item(62, 13).Value = -4375.5
$ws.Cells.Item(62, 14).Value = -13248
$ws.Cells.Item(65, 8).Value = 9999.857
$ws.Cells.Item(65, 9).Value = 4999.5
$ws.Cells.Item(65, 10).Value = 12000
$ws.Cells.Item(65, 11).Value = 24997.5
$ws.Cells.Item(65, 12).Value = 60000
$ws.Cells.Item(65, 13).Value = -21877.5
$ws.Cells.Item(65, 14).Value = -66240
$ws.Cells.Item(70, 8).Value = 41249.75
$ws.Cells.Item(70, 10).Value = 40000
$ws.Cells.Item(70, 12).Value = 40000
$ws.Cells.Item(70, 14).Value = -40630
$ws.Cells.Item(73, 8).Value = 41249.75
$ws.Cells.Item(73, 10).Value = 40000
$ws.Cells.Item(73, 12).Value = 40000
$ws.Cells.Item(73, 14).Value = -42184
$ws.Cells.Item(81, 8).Value = 678.6
$ws.Cells.Item(81, 9).Value = 678.6
$ws.Cells.Item(81, 11).Value = 1357.2
$ws.Cells.Item(81, 13).Value = -296.2
$ws.Cells.Item(84, 8).Value = 678.6
$ws.Cells.Item(84, 9).Value = 678.6
$ws.Cells.Item(84, 11).Value = 6786
$ws.Cells.Item(84, 13).Value = -1482
$ws.Cells.Item(126, 8).Value = 4383
$ws.Cells.Item(126, 9).Value = 1799.4286
$ws.Cells.Item(126, 11).Value = 5398.2858
$ws.Cells.Item(126, 13).Value = -2928.2858
$ws.Cells.Item(132, 8).Value = 1364.3636
$ws.Cells.Item(132, 9).Value = 1000.8889
$ws.Cells.Item(132, 11).Value = 3002.6667
$ws.Cells.Item(132, 13).Value = -472.6667000000002
$ws.Cells.Item(135, 8).Value = 71810
$ws.Cells.Item(135, 10).Value = 71810
$ws.Cells.Item(135, 12).Value = 71810
$ws.Cells.Item(135, 14).Value = -81950
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).Value = $null
